$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 6.8
$ws.Cells.Item(2, 8).Value = 1.62
$ws.Cells.Item(2, 9).Value = 1.82
$ws.Cells.Item(2, 10).Value = 3.6
$ws.Cells.Item(2, 11).Value = 4.5
$ws.Cells.Item(2, 12).Value = 1.39
$ws.Cells.Item(2, 13).Value = 1.07
$ws.Cells.Item(2, 14).Value = 3.2
$ws.Cells.Item(2, 15).Value = 1.32
$ws.Cells.Item(2, 16).Value = 1.83
$ws.Cells.Item(2, 17).Value = 1.94
$ws.Cells.Item(2, 18).Value = 1.33
$ws.Cells.Item(2, 19).Value = 3.1
$ws.Cells.Item(2, 20).Value = 1.91
$ws.Cells.Item(2, 21).Value = 1.85
$ws.Cells.Item(2, 22).Value = 2.2
$ws.Cells.Item(2, 23).Value = 1.17
$ws.Cells.Item(2, 24).Value = 15
$ws.Cells.Item(2, 25).Value = 8.6
$ws.Cells.Item(2, 27).Value = 18.5
$ws.Cells.Item(2, 28).Value = 19.5
$ws.Cells.Item(2, 29).Value = 9.6
$ws.Cells.Item(2, 30).Value = 11
$ws.Cells.Item(2, 32).Value = 55
$ws.Cells.Item(2, 33).Value = 25
$ws.Cells.Item(2, 34).Value = 24
$ws.Cells.Item(2, 35).Value = 44
$ws.Cells.Item(2, 37).Value = 260
$ws.Cells.Item(2, 38).Value = 110
$ws.Cells.Item(2, 41).Value = 12

# Row 3
$ws.Cells.Item(3, 8).Value = 1.04
$ws.Cells.Item(3, 10).Value = 3.05
$ws.Cells.Item(3, 11).Value = 950
$ws.Cells.Item(3, 12).Value = 1.26
$ws.Cells.Item(3, 13).Value = 1.03
$ws.Cells.Item(3, 18).Value = 1.27

# Row 4
$ws.Cells.Item(4, 6).Value = 5.3
$ws.Cells.Item(4, 7).Value = 7
$ws.Cells.Item(4, 8).Value = 1.6
$ws.Cells.Item(4, 9).Value = 1.79
$ws.Cells.Item(4, 10).Value = 3.65
$ws.Cells.Item(4, 11).Value = 4.5
$ws.Cells.Item(4, 22).Value = 2.26
$ws.Cells.Item(4, 23).Value = 1.16
$ws.Cells.Item(4, 24).Value = 16
$ws.Cells.Item(4, 25).Value = 8.6
$ws.Cells.Item(4, 26).Value = 10.5
$ws.Cells.Item(4, 27).Value = 18
$ws.Cells.Item(4, 28).Value = 21
$ws.Cells.Item(4, 30).Value = 11
$ws.Cells.Item(4, 31).Value = 20
$ws.Cells.Item(4, 32).Value = 55
$ws.Cells.Item(4, 33).Value = 26
$ws.Cells.Item(4, 34).Value = 24
$ws.Cells.Item(4, 35).Value = 42
$ws.Cells.Item(4, 37).Value = 250
$ws.Cells.Item(4, 38).Value = 330
$ws.Cells.Item(4, 41).Value = 11.5

# Row 5
$ws.Cells.Item(5, 6).Value = 1.36
$ws.Cells.Item(5, 7).Value = 2.3
$ws.Cells.Item(5, 8).Value = 2.84
$ws.Cells.Item(5, 9).Value = 44
$ws.Cells.Item(5, 10).Value = 2.24
$ws.Cells.Item(5, 11).Value = 44
$ws.Cells.Item(5, 13).Value = 1.02
$ws.Cells.Item(5, 15).Value = 1.03
$ws.Cells.Item(5, 17).Value = 1.36
$ws.Cells.Item(5, 22).Value = 1.1
$ws.Cells.Item(5, 23).Value = 1.76
$ws.Cells.Item(5, 29).Value = 42

# Row 6
$ws.Cells.Item(6, 6).Value = 2.12
$ws.Cells.Item(6, 7).Value = 2.58
$ws.Cells.Item(6, 8).Value = 2.88
$ws.Cells.Item(6, 10).Value = 3.15
$ws.Cells.Item(6, 11).Value = 4.1
$ws.Cells.Item(6, 17).Value = 1.82
$ws.Cells.Item(6, 20).Value = 1.89
$ws.Cells.Item(6, 22).Value = 1.39
$ws.Cells.Item(6, 23).Value = 1.63
$ws.Cells.Item(6, 24).Value = 15.5
$ws.Cells.Item(6, 25).Value = 14.5
$ws.Cells.Item(6, 26).Value = 27
$ws.Cells.Item(6, 27).Value = 290
$ws.Cells.Item(6, 28).Value = 11
$ws.Cells.Item(6, 29).Value = 8.800000000000001
$ws.Cells.Item(6, 30).Value = 16.5
$ws.Cells.Item(6, 31).Value = 120
$ws.Cells.Item(6, 32).Value = 16.5
$ws.Cells.Item(6, 33).Value = 12.5
$ws.Cells.Item(6, 34).Value = 20
$ws.Cells.Item(6, 35).Value = 150
$ws.Cells.Item(6, 36).Value = 36
$ws.Cells.Item(6, 37).Value = 29
$ws.Cells.Item(6, 38).Value = 110
$ws.Cells.Item(6, 39).Value = 580
$ws.Cells.Item(6, 40).Value = 22
$ws.Cells.Item(6, 41).Value = 140

# Row 7
$ws.Cells.Item(7, 8).Value = 1.01
$ws.Cells.Item(7, 11).Value = 1000
$ws.Cells.Item(7, 23).Value = 1.01

# Row 8
$ws.Cells.Item(8, 7).Value = 2.06
$ws.Cells.Item(8, 10).Value = 3.1
$ws.Cells.Item(8, 11).Value = 3.9
$ws.Cells.Item(8, 16).Value = 1.65
$ws.Cells.Item(8, 22).Value = 1.2
$ws.Cells.Item(8, 23).Value = 1.94
$ws.Cells.Item(8, 28).Value = 7.2
$ws.Cells.Item(8, 29).Value = 8.800000000000001
$ws.Cells.Item(8, 32).Value = 11.5
$ws.Cells.Item(8, 33).Value = 12
$ws.Cells.Item(8, 36).Value = 25
$ws.Cells.Item(8, 37).Value = 28
$ws.Cells.Item(8, 38).Value = 160
$ws.Cells.Item(8, 40).Value = 24

# Row 9
$ws.Cells.Item(9, 6).Value = 2.36
$ws.Cells.Item(9, 7).Value = 2.64
$ws.Cells.Item(9, 8).Value = 2.82
$ws.Cells.Item(9, 9).Value = 3.35
$ws.Cells.Item(9, 14).Value = 3.55
$ws.Cells.Item(9, 15).Value = 1.28
$ws.Cells.Item(9, 18).Value = 1.34
$ws.Cells.Item(9, 20).Value = 1.68
$ws.Cells.Item(9, 21).Value = 2.08
$ws.Cells.Item(9, 23).Value = 1.6
$ws.Cells.Item(9, 24).Value = 90
$ws.Cells.Item(9, 29).Value = 14
$ws.Cells.Item(9, 34).Value = 60
$ws.Cells.Item(9, 36).Value = 170
$ws.Cells.Item(9, 39).Value = 580

# Row 10
$ws.Cells.Item(10, 6).Value = 1.09
$ws.Cells.Item(10, 22).Value = 1.41

# Row 11
$ws.Cells.Item(11, 6).Value = 1.61
$ws.Cells.Item(11, 16).Value = 2.86
$ws.Cells.Item(11, 18).Value = 1.75
$ws.Cells.Item(11, 19).Value = 2.2
$ws.Cells.Item(11, 21).Value = 2.48
$ws.Cells.Item(11, 22).Value = 1.2
$ws.Cells.Item(11, 24).Value = 28
$ws.Cells.Item(11, 28).Value = 26
$ws.Cells.Item(11, 32).Value = 13.5
$ws.Cells.Item(11, 33).Value = 11
$ws.Cells.Item(11, 34).Value = 25
$ws.Cells.Item(11, 35).Value = 130
$ws.Cells.Item(11, 36).Value = 17
$ws.Cells.Item(11, 38).Value = 26
$ws.Cells.Item(11, 39).Value = 580
$ws.Cells.Item(11, 41).Value = 980

# Row 12
$ws.Cells.Item(12, 6).Value = 2.74
$ws.Cells.Item(12, 7).Value = 2.92
$ws.Cells.Item(12, 8).Value = 2.66
$ws.Cells.Item(12, 9).Value = 2.82
$ws.Cells.Item(12, 10).Value = 3.45
$ws.Cells.Item(12, 15).Value = 1.24
$ws.Cells.Item(12, 16).Value = 2.2
$ws.Cells.Item(12, 17).Value = 1.72
$ws.Cells.Item(12, 18).Value = 1.48
$ws.Cells.Item(12, 20).Value = 1.6
$ws.Cells.Item(12, 21).Value = 2.44
$ws.Cells.Item(12, 22).Value = 1.55
$ws.Cells.Item(12, 23).Value = 1.52
$ws.Cells.Item(12, 25).Value = 15
$ws.Cells.Item(12, 26).Value = 24
$ws.Cells.Item(12, 27).Value = 280
$ws.Cells.Item(12, 29).Value = 8.800000000000001
$ws.Cells.Item(12, 33).Value = 13
$ws.Cells.Item(12, 34).Value = 15
$ws.Cells.Item(12, 35).Value = 75
$ws.Cells.Item(12, 36).Value = 120
$ws.Cells.Item(12, 37).Value = 32
$ws.Cells.Item(12, 38).Value = 95
$ws.Cells.Item(12, 39).Value = 580
$ws.Cells.Item(12, 40).Value = 21

# Row 13
$ws.Cells.Item(13, 6).Value = 3.5
$ws.Cells.Item(13, 7).Value = 3.75
$ws.Cells.Item(13, 9).Value = 2.16
$ws.Cells.Item(13, 10).Value = 4
$ws.Cells.Item(13, 13).Value = 1.04
$ws.Cells.Item(13, 15).Value = 1.2
$ws.Cells.Item(13, 16).Value = 2.52
$ws.Cells.Item(13, 17).Value = 1.58
$ws.Cells.Item(13, 19).Value = 2.38
$ws.Cells.Item(13, 22).Value = 1.87
$ws.Cells.Item(13, 24).Value = 24
$ws.Cells.Item(13, 26).Value = 16.5
$ws.Cells.Item(13, 27).Value = 27
$ws.Cells.Item(13, 28).Value = 23
$ws.Cells.Item(13, 29).Value = 9.800000000000001
$ws.Cells.Item(13, 30).Value = 11
$ws.Cells.Item(13, 31).Value = 19
$ws.Cells.Item(13, 32).Value = 30
$ws.Cells.Item(13, 33).Value = 15.5
$ws.Cells.Item(13, 34).Value = 15
$ws.Cells.Item(13, 35).Value = 27
$ws.Cells.Item(13, 36).Value = 65
$ws.Cells.Item(13, 38).Value = 38
$ws.Cells.Item(13, 39).Value = 60
$ws.Cells.Item(13, 40).Value = 28

# Row 14
$ws.Cells.Item(14, 6).Value = 5.6
$ws.Cells.Item(14, 7).Value = 6.2
$ws.Cells.Item(14, 8).Value = 1.64
$ws.Cells.Item(14, 9).Value = 1.68
$ws.Cells.Item(14, 10).Value = 4.3
$ws.Cells.Item(14, 14).Value = 5.3
$ws.Cells.Item(14, 15).Value = 1.2
$ws.Cells.Item(14, 16).Value = 2.46
$ws.Cells.Item(14, 17).Value = 1.61
$ws.Cells.Item(14, 18).Value = 1.61
$ws.Cells.Item(14, 19).Value = 2.48
$ws.Cells.Item(14, 20).Value = 1.69
$ws.Cells.Item(14, 21).Value = 2.32
$ws.Cells.Item(14, 22).Value = 2.46
$ws.Cells.Item(14, 23).Value = 1.2
$ws.Cells.Item(14, 25).Value = 12
$ws.Cells.Item(14, 26).Value = 11.5
$ws.Cells.Item(14, 27).Value = 17
$ws.Cells.Item(14, 28).Value = 27
$ws.Cells.Item(14, 31).Value = 15
$ws.Cells.Item(14, 32).Value = 55
$ws.Cells.Item(14, 33).Value = 22
$ws.Cells.Item(14, 34).Value = 18
$ws.Cells.Item(14, 36).Value = 140
$ws.Cells.Item(14, 37).Value = 70
$ws.Cells.Item(14, 38).Value = 150
$ws.Cells.Item(14, 39).Value = 80
$ws.Cells.Item(14, 40).Value = 65
$ws.Cells.Item(14, 41).Value = 7

# Row 15
$ws.Cells.Item(15, 14).Value = 1.02
